# Update the confidential-notice footer text: roll the "as of" date
# forward one day (2021-03-24 -> 2021-03-25). The sheet is protected, so
# we briefly unprotect it to make the edits, then restore protection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$oldText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-03-24 for illustrative purposes only and are subject to change."
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-03-25 for illustrative purposes only and are subject to change."

for ($r = 16; $r -le 35; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq $oldText) {
        $cell.Value = $newText
    }
}

# Refreshed Weight / Percent Change figures for rows 2-13 (Symbol table).
$ws.Range("D2").Value = 0.03042094874181618
$ws.Range("E2").Value = 0.01040118870728057

$ws.Range("D3").Value = 0.02421087903042209
$ws.Range("E3").Value = 0.01084545230465861

$ws.Range("D4").Value = 0.05260751442686715
$ws.Range("E4").Value = 0.01033973412112266

$ws.Range("D5").Value = 0.1383706988420288
$ws.Range("E5").Value = 0.004856052722858051

$ws.Range("D6").Value = 0.03138923661560406
$ws.Range("E6").Value = 0.003782148260211615

$ws.Range("D7").Value = 0.119166957538177
$ws.Range("E7").Value = 0.01258457374830835

$ws.Range("D8").Value = 0.1008695418103936
$ws.Range("E8").Value = 0.01491646778042965

$ws.Range("D9").Value = 0.02766362649330089
$ws.Range("E9").Value = 0.01661985018726586

$ws.Range("D10").Value = 0.1215121289293207
$ws.Range("E10").Value = 0.01713688610240327

$ws.Range("D11").Value = 0.2489665437874733
$ws.Range("E11").Value = -0.00009689922480626834

$ws.Range("D12").Value = 0.1048219237845962
$ws.Range("E12").Value = -0.002849002849002913

$ws.Range("E13").Value = 0.007137218682597535

$ws.Protect()
